$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Remove the row for "NURIA BARBOSA MARCAL" (CALDAS NOVAS / GO551004929),
# shifting all rows below it up by one, as seen in the source diff.
$ws.Rows.Item(520).Delete()
